# Apply updated cryptocurrency market data (price/volume) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new "Price" values are plain decimal numbers (e.g. "503.58") that Excel
# would otherwise auto-convert to a Number when assigned via .Value. The sheet
# stores these as text, so first force Text format on just those cells (reusing
# a single Union range to avoid creating one new style per cell).
$priceCell0 = $ws.Range("D5")
$priceCell1 = $ws.Range("D6")
$priceCell2 = $ws.Range("D7")
$priceCell3 = $ws.Range("D9")
$priceCell4 = $ws.Range("D11")
$priceCell5 = $ws.Range("D14")
$priceCell6 = $ws.Range("D18")
$priceCell7 = $ws.Range("D20")
$priceCell8 = $ws.Range("D21")
$priceCell9 = $ws.Range("D23")
$priceCell10 = $ws.Range("D26")
$priceCell11 = $ws.Range("D27")
$priceCell12 = $ws.Range("D30")
$priceCell13 = $ws.Range("D31")
$priceCell14 = $ws.Range("D33")
$priceCell15 = $ws.Range("D34")
$priceCell16 = $ws.Range("D41")
$priceCell17 = $ws.Range("D42")
$priceCell18 = $ws.Range("D43")
$priceCell19 = $ws.Range("D45")
$priceCell20 = $ws.Range("D46")
$priceCell21 = $ws.Range("D50")
$priceCell22 = $ws.Range("D51")
$textPriceRange = $excel.Union($priceCell0, $priceCell1, $priceCell2, $priceCell3, $priceCell4, $priceCell5, $priceCell6, $priceCell7, $priceCell8, $priceCell9, $priceCell10, $priceCell11, $priceCell12, $priceCell13, $priceCell14, $priceCell15, $priceCell16, $priceCell17, $priceCell18, $priceCell19, $priceCell20, $priceCell21, $priceCell22)
foreach ($area in $textPriceRange.Areas) {
    $area.NumberFormat = "@"
}

# Row 51 previously listed "dogwifhat"; it is now replaced by "ZEEBU".
$ws.Range("B51").Value = "ZEEBU"
$ws.Range("C51").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"

# Update Price (D) and Volume(1h) (E) columns row by row.
$ws.Range("D2").Value = '54.466.79'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '2.289.47'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '503.58'
$ws.Range("E5").Value = '  +1.80%  '
$ws.Range("D6").Value = '130.32'
$ws.Range("E6").Value = '  +2.26%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("D9").Value = '0.0959'
$ws.Range("E9").Value = '  +1.78%  '
$ws.Range("E10").Value = '  +0.78%  '
$ws.Range("D11").Value = '0.339'
$ws.Range("E11").Value = '  +4.95%  '
$ws.Range("E12").Value = '  +2.52%  '
$ws.Range("D13").Value = '2.700.07'
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").Value = '22.85'
$ws.Range("E14").Value = '  +6.09%  '
$ws.Range("D15").Value = '54.452.93'
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("D17").Value = '2.293.74'
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("D18").Value = '10.25'
$ws.Range("E18").Value = '  +2.98%  '
$ws.Range("E19").Value = '  +2.82%  '
$ws.Range("D20").Value = '303.73'
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").Value = '6.38'
$ws.Range("E21").Value = '  -0.61%  '
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").Value = '61.94'
$ws.Range("E23").Value = '  -2.73%  '
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("E25").Value = '  +2.12%  '
$ws.Range("D26").Value = '7.35'
$ws.Range("E26").Value = '  +3.56%  '
$ws.Range("D27").Value = '171.68'
$ws.Range("E27").Value = '  +2.33%  '
$ws.Range("E28").Value = '  +2.15%  '
$ws.Range("E29").Value = '  +1.40%  '
$ws.Range("D30").Value = '5.97'
$ws.Range("E30").Value = '  +1.68%  '
$ws.Range("D31").Value = '1.09'
$ws.Range("E31").Value = '  +1.40%  '
$ws.Range("D33").Value = '17.86'
$ws.Range("E33").Value = '  +1.67%  '
$ws.Range("D34").Value = '0.966'
$ws.Range("E34").Value = '  +10.92%  '
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("E36").Value = '  +0.88%  '
$ws.Range("E37").Value = '  +3.27%  '
$ws.Range("E38").Value = '  +0.63%  '
$ws.Range("E39").Value = '  +1.27%  '
$ws.Range("E40").Value = '  +1.56%  '
$ws.Range("D41").Value = '4.95'
$ws.Range("E41").Value = '  +2.85%  '
$ws.Range("D42").Value = '126.22'
$ws.Range("D43").Value = '0.0496'
$ws.Range("E43").Value = '  +3.43%  '
$ws.Range("E44").Value = '  +1.04%  '
$ws.Range("D45").Value = '0.549'
$ws.Range("E45").Value = '  +1.15%  '
$ws.Range("D46").Value = '242.86'
$ws.Range("E46").Value = '  +1.75%  '
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("E48").Value = '  +1.54%  '
$ws.Range("E49").Value = '  +0.86%  '
$ws.Range("D50").Value = '16.46'
$ws.Range("E50").Value = '  +1.28%  '
$ws.Range("D51").Value = '4.64'
$ws.Range("E51").Value = '  -0.48%  '
